{"js": "// Replace the multiplication-problem text in each table cell with the\n// newly generated values, preserving all existing run/paragraph\n// formatting (font, size, justification, etc.) by doing an in-place\n// text replacement rather than rebuilding the paragraphs.\nconst replacements = [\n  [\"99\u00d788=8712\", \"37\u00d714=518\"],\n  [\"15\u00d745=675\", \"89\u00d752=4628\"],\n  [\"89\u00d727=2403\", \"71\u00d797=6887\"],\n  [\"13\u00d717=221\", \"98\u00d752=5096\"],\n  [\"91\u00d752=4732\", \"50\u00d752=2600\"],\n  [\"66\u00d774=4884\", \"98\u00d793=9114\"],\n  [\"99\u00d783=8217\", \"25\u00d743=1075\"],\n  [\"96\u00d788=8448\", \"28\u00d720=560\"],\n  [\"18\u00d768=1224\", \"98\u00d784=8232\"],\n  [\"34\u00d747=1598\", \"35\u00d772=2520\"],\n  [\"93\u00d764=5952\", \"48\u00d751=2448\"],\n  [\"31\u00d718=558\", \"98\u00d795=9310\"],\n  [\"84\u00d726=2184\", \"52\u00d741=2132\"],\n  [\"78\u00d734=2652\", \"33\u00d745=1485\"],\n  [\"83\u00d740=3320\", \"98\u00d763=6174\"],\n  [\"18\u00d751=918\", \"90\u00d763=5670\"],\n  [\"88\u00d767=5896\", \"94\u00d734=3196\"],\n  [\"41\u00d726=1066\", \"54\u00d799=5346\"],\n  [\"29\u00d745=1305\", \"53\u00d750=2650\"],\n  [\"28\u00d745=1260\", \"87\u00d767=5829\"],\n  [\"35\u00d745=1575\", \"78\u00d757=4446\"],\n  [\"73\u00d719=1387\", \"91\u00d754=4914\"],\n  [\"29\u00d713=377\", \"65\u00d778=5070\"],\n  [\"40\u00d778=3120\", \"80\u00d751=4080\"],\n  [\"37\u00d766=2442\", \"76\u00d755=4180\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the multiplication-problem text in each table cell with the\n# newly generated values, preserving all existing run/paragraph\n# formatting (font, size, justification, etc.) via Find/Replace on the\n# document's content range (mirrors Word's Ctrl+H \"Replace All\").\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"99\u00d788=8712\", \"37\u00d714=518\"),\n    @(\"15\u00d745=675\", \"89\u00d752=4628\"),\n    @(\"89\u00d727=2403\", \"71\u00d797=6887\"),\n    @(\"13\u00d717=221\", \"98\u00d752=5096\"),\n    @(\"91\u00d752=4732\", \"50\u00d752=2600\"),\n    @(\"66\u00d774=4884\", \"98\u00d793=9114\"),\n    @(\"99\u00d783=8217\", \"25\u00d743=1075\"),\n    @(\"96\u00d788=8448\", \"28\u00d720=560\"),\n    @(\"18\u00d768=1224\", \"98\u00d784=8232\"),\n    @(\"34\u00d747=1598\", \"35\u00d772=2520\"),\n    @(\"93\u00d764=5952\", \"48\u00d751=2448\"),\n    @(\"31\u00d718=558\", \"98\u00d795=9310\"),\n    @(\"84\u00d726=2184\", \"52\u00d741=2132\"),\n    @(\"78\u00d734=2652\", \"33\u00d745=1485\"),\n    @(\"83\u00d740=3320\", \"98\u00d763=6174\"),\n    @(\"18\u00d751=918\", \"90\u00d763=5670\"),\n    @(\"88\u00d767=5896\", \"94\u00d734=3196\"),\n    @(\"41\u00d726=1066\", \"54\u00d799=5346\"),\n    @(\"29\u00d745=1305\", \"53\u00d750=2650\"),\n    @(\"28\u00d745=1260\", \"87\u00d767=5829\"),\n    @(\"35\u00d745=1575\", \"78\u00d757=4446\"),\n    @(\"73\u00d719=1387\", \"91\u00d754=4914\"),\n    @(\"29\u00d713=377\", \"65\u00d778=5070\"),\n    @(\"40\u00d778=3120\", \"80\u00d751=4080\"),\n    @(\"37\u00d766=2442\", \"76\u00d755=4180\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute(\n        $oldText,\n        $false, $false, $false, $false, $false, $true, 1, $false,\n        $newText, 2\n    )\n}\n"}
